# Update the "About" sheet's source-link cell (B6) so its displayed text
# points at the new CEPE working-paper URL. The cell keeps its existing
# hyperlink (same relationship/target) — only the visible text changes,
# which is exactly what the authoring diff shows: the old URL string is
# dropped from the shared-string table and the new URL string is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("B6").Value = "https://ethz.ch/content/dam/ethz/special-interest/mtec/cepe/cepe-dam/documents/research/cepe-wp/CEPE_WP86.pdf"
